# Error Calculations and Plots
# Re-derive which cells are missing (simulating a different random seed
# of removed data) and drop two rows ("RM 232" and "SC 92") from the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply scattered cell value / empty changes (before row deletions) ---
$ws.Range("E2").Value = -7.2
$ws.Range("F3").Value = $null
$ws.Range("F4").Value = 17.97
$ws.Range("F5").Value = $null
$ws.Range("E6").Value = $null
$ws.Range("F8").Value = $null
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = $null
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E22").Value = $null
$ws.Range("E23").Value = $null
$ws.Range("F23").Value = 16.48

$ws.Range("F29").Value = $null      # SC 101 F: 17 -> empty
$ws.Range("F31").Value = 18.06      # SC 119 F: empty -> 18.06
$ws.Range("B32").Value = -19.7      # SC 120 B: empty -> -19.7
$ws.Range("E33").Value = -8.1       # SC 132 E: empty -> -8.1
$ws.Range("B34").Value = $null      # SC 193 B: -19.9 -> empty
$ws.Range("E35").Value = -10.7      # SC 232 E: empty -> -10.7

# --- Delete the two rows that were removed entirely (RM 232, SC 92) ---
$ws.Rows.Item(28).EntireRow.Delete()  # SC 92  (row 28, after RM232 still above it)
$ws.Rows.Item(26).EntireRow.Delete()  # RM 232 (row 26)
